$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.127.62'
$ws.Range("E2").Value = '  +1.14%  '
$ws.Range("D3").Value = '3.403.33'
$ws.Range("E3").Value = '  +1.39%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = "'581.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.55%  '
$ws.Range("D6").Value = "'178.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.19%  '
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("E8").Value = '  +0.43%  '
$ws.Range("D9").Value = "'0.198"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.43%  '
$ws.Range("D10").Value = "'0.585"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.63%  '
$ws.Range("D11").Value = "'48.35"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.67%  '
$ws.Range("D12").Value = "'0.0000281"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.86%  '
$ws.Range("D13").Value = "'680.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.33%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '3.955.37'
$ws.Range("E14").Value = '  +1.26%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = "'8.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.79%  '
$ws.Range("D16").Value = '69.354.42'
$ws.Range("E16").Value = '  +1.44%  '
$ws.Range("D17").Value = '3.408.46'
$ws.Range("E17").Value = '  +1.56%  '
$ws.Range("E18").Value = '  +0.48%  '
$ws.Range("D20").Value = "'11.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.69%  '
$ws.Range("D21").Value = "'0.910"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.61%  '
$ws.Range("D22").Value = "'5.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.22%  '
$ws.Range("D23").Value = "'17.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.42%  '
$ws.Range("D24").Value = "'100.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.45%  '
$ws.Range("E25").Value = '  -0.46%  '
$ws.Range("D26").Value = "'2.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.35%  '
$ws.Range("D27").Value = "'9.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.74%  '
$ws.Range("D28").Value = "'33.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.31%  '
$ws.Range("D29").Value = "'8.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.38%  '
$ws.Range("D30").Value = "'6.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.57%  '
$ws.Range("D31").Value = "'3.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +10.78%  '
$ws.Range("D32").Value = "'557.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.91%  '
$ws.Range("D33").Value = "'11.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.77%  '
$ws.Range("E34").Value = '  -0.40%  '
$ws.Range("D35").Value = "'58.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("E36").Value = '  +0.12%  '
$ws.Range("D37").Value = '3.608.38'
$ws.Range("E37").Value = '  -2.99%  '
$ws.Range("E38").Value = '  +0.20%  '
$ws.Range("D39").Value = "'34.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.55%  '
$ws.Range("D40").Value = '0.0₃0739'
$ws.Range("E40").Value = '  +9.74%  '
$ws.Range("D41").Value = "'3.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.97%  '
$ws.Range("D42").Value = "'2.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.80%  '
$ws.Range("E43").Value = '  +2.59%  '
$ws.Range("D45").Value = "'0.335"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.13%  '
$ws.Range("D46").Value = "'2.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.73%  '
$ws.Range("E47").Value = '  +0.18%  '
$ws.Range("E48").Value = '  +3.78%  '
$ws.Range("E49").Value = '  -0.12%  '
$ws.Range("D50").Value = "'131.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.37%  '
$ws.Range("D51").Value = "'2.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.58%  '
